# Updated symbol list on Fri Jan 20 03:37:52 UTC 2023 with GitHub Actions
# Refresh of crypto price/volume columns (D = Price, E = Volume(1h)).
# Values are entered with a leading apostrophe to force text storage,
# matching the workbook's existing inlineStr (text) cell type so that
# formatting such as trailing zeros (e.g. "0.08300") is preserved exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'295.11"
$ws.Range("E2").Value = "'1.77%"
$ws.Range("D3").Value = "'30.97"
$ws.Range("E3").Value = "'0.55%"
$ws.Range("D4").Value = "'4.909"
$ws.Range("E4").Value = "'-0.68%"
$ws.Range("D5").Value = "'0.07445"
$ws.Range("E5").Value = "'4.42%"
$ws.Range("D6").Value = "'2.155"
$ws.Range("E6").Value = "'19.61%"
$ws.Range("D7").Value = "'7.744"
$ws.Range("E7").Value = "'0.80%"
$ws.Range("D8").Value = "'3.750"
$ws.Range("E8").Value = "'0.44%"
$ws.Range("E9").Value = "'1.92%"
$ws.Range("D10").Value = "'0.08900"
$ws.Range("E10").Value = "'16.97%"
$ws.Range("D11").Value = "'0.1713"
$ws.Range("D12").Value = "'0.08319"
$ws.Range("D13").Value = "'0.03150"
$ws.Range("E13").Value = "'3.30%"
$ws.Range("E14").Value = "'0.56%"
$ws.Range("D15").Value = "'0.001519"
$ws.Range("E15").Value = "'1.17%"
$ws.Range("D16").Value = "'0.005716"
$ws.Range("E16").Value = "'-0.09%"
$ws.Range("E17").Value = "'1.16%"
$ws.Range("E18").Value = "'-0.31%"
$ws.Range("D20").Value = "'0.1293"
$ws.Range("E20").Value = "'-0.23%"
$ws.Range("D21").Value = "'3.965"
$ws.Range("E21").Value = "'-1.87%"
$ws.Range("E23").Value = "'0.76%"
$ws.Range("E24").Value = "'0.23%"
$ws.Range("D25").Value = "'0.004621"
$ws.Range("E25").Value = "'15.50%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'4.06%"
$ws.Range("D27").Value = "'0.0003398"
$ws.Range("E27").Value = "'-95.49%"
$ws.Range("D39").Value = "'0.01613"
$ws.Range("E39").Value = "'-0.04%"
$ws.Range("D40").Value = "'0.04487"
$ws.Range("E40").Value = "'2.77%"
$ws.Range("D41").Value = "'0.007259"
$ws.Range("E41").Value = "'-1.18%"
$ws.Range("D42").Value = "'0.008970"
$ws.Range("D43").Value = "'0.1333"
$ws.Range("E43").Value = "'2.28%"
$ws.Range("D44").Value = "'0.001964"
$ws.Range("E44").Value = "'-1.78%"
$ws.Range("E45").Value = "'-0.15%"
$ws.Range("D46").Value = "'0.00006114"
$ws.Range("E46").Value = "'2.05%"
$ws.Range("E47").Value = "'0.10%"
$ws.Range("D48").Value = "'2.236"
$ws.Range("E48").Value = "'-0.43%"
$ws.Range("D49").Value = "'0.002002"
$ws.Range("E49").Value = "'-33.23%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.10%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.10%"
